$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 227, shifting the existing rows 227:397 down to 228:398.
$ws.Rows("227").Insert()

# Populate the newly inserted row 227 with a new record (same series as its
# neighbours, but a new "Fecha" of 44978, i.e. 2023-02-21) using the same
# Volumen/Precio values that previously belonged to row 227 before the shift.
$ws.Range("A227").Value = 3
$ws.Range("B227").Value = "Femacal de La Calera"
$ws.Range("C227").Value = "Coquimbo"
$ws.Range("D227").Value = 44978
$ws.Range("E227").Value = 5
$ws.Range("F227").Value = 100112039
$ws.Range("G227").Value = "Ciboulette"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 120
$ws.Range("K227").Value = 1500
$ws.Range("L227").Value = 1500
$ws.Range("M227").Value = 1500
$ws.Range("N227").Value = "$/docena de atados"
$ws.Range("O227").Value = "Provincia de Quillota"
$ws.Range("P227").Value = 500
$ws.Range("Q227").Value = 3
$ws.Range("R227").Value = "Hortaliza"
